$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027941584587097
$ws.Range("B1").Value = 1.732664465904236
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.051704883575439
